$d = $word.ActiveDocument

# --- Split "Version" into "Versi" + "on" runs ---
# A temporary bookmark at the split point stops the two halves from being
# re-merged into a single run; once the split exists we can drop the marker.
$d.Bookmarks.Add("tmpSplit", $d.Range(5, 5))
$d.Bookmarks("tmpSplit").Delete()

# --- "Version 1." -> "Version 2." ---
$d.Content.Find.Execute("1", $false, $false, $false, $false, $false, $true, 1, $false, "2", 2) | Out-Null

# --- Move the trailing "." to after the _GoBack bookmark ---
$r = $d.Range(9, 10)
$r.Text = ""
$r2 = $d.Range(10, 10)
$r2.InsertAfter(".")
